$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Career-summary sentence tweak.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Recognized for leadership, mentorship, and cross-functional collaboration, with a strategic approach to system design and a strong commitment to continuous learning in cybersecurity.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Recognized for leadership, mentorship, cross-functional collaboration, and a strong commitment to continuous learning in cybersecurity.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a new "in progress" OSCP certification entry at the top of the
#    EDUCATION & CERTIFICATIONS list (right before "Enterprise Penetration
#    Testing (GPEN) ...").
# ---------------------------------------------------------------------------
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Enterprise Penetration Testing (*") {
        $idx = $i
        break
    }
}

# Push everything down one slot, leaving a fresh empty paragraph at $idx.
$d.Paragraphs.Item($idx).Range.InsertParagraphBefore()

$d.Paragraphs.Item($idx).Range.InsertBefore("Offensive Security Certified Professional (")
$d.Paragraphs.Item($idx).Range.InsertAfter("OSCP")
$d.Paragraphs.Item($idx).Range.InsertAfter(") | ")
$d.Paragraphs.Item($idx).Range.InsertAfter("OffSec")
$d.Paragraphs.Item($idx).Range.InsertAfter("`t")
$d.Paragraphs.Item($idx).Range.InsertAfter("Expected Aug. 2025")

# Bold just the "OSCP" run (mirrors the bold GPEN/GCIH/GSEC/GFACT runs).
$fr = $d.Paragraphs.Item($idx).Range.Duplicate
$fr.Find.Execute("OSCP", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Bold = 1
$fr.BoldBi = 1
